$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Update the product name (shared string referenced by both sheets' B1)
$newProductName = "4274-MS-EI-DB-DL-REC-RNI-FEE+INT-FFConMONTHLYonLASTSUN-FIFC-1-FFROP-DL-FIFR-1-MD-TR-1-1st"
$ws1.Range("B1").Value = $newProductName
$ws2.Range("B1").Value = $newProductName

# Update the short name from the numeric 4274 to the text "427m"
$ws1.Range("B2").Value = "427m"

# Make ProductLoanInput the active/selected sheet (was ProductLoanOutput before)
$ws1.Activate()
$ws1.Range("B3").Select()
